$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# uv2vxvy (row 7) status is now "Done"
$ws.Range("B7").Value = "Done"

# mapzoomps (row 22) status + notes
$ws.Range("B22").Value = "Started"
$ws.Range("C22").Value = "Needs scarlabel and scalebarps"

# New column D notes describing which functions/tools are needed for testing
$ws.Range("D14").Value = "needs mapzoomps to test"
$ws.Range("D11").Value = "bedmap2 for testing"
$ws.Range("D5").Value  = "antmap, bedmap2, pcolorm for testing"
$ws.Range("D12").Value = "pcolorps, graticuleps for testing"
$ws.Range("D13").Value = "plotps for testing"
$ws.Range("D17").Value = "needs patch, uistack for testing"
$ws.Range("D7").Value  = "needs antmap, pcolorm, bedmap2, quivermc, quiver for testing"

# Widen column D to fit the new notes (matches Excel's auto-fit width of 54)
$ws.Columns.Item(4).ColumnWidth = 53.14

# Update view: scroll back to top and select B7
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()
